$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Escopo")

$r = $ws.Range("A2:A5")
$r.MergeCells = $false
$r.MergeCells = $true
$r.Borders.LineStyle = 1
$r.Borders.Weight = 2
$r.Borders.ColorIndex = -4105
